$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" metadata value (row 8, column B)
$ws.Range("B8").Value = "2025-01-15T07:47:50+00:00"

# Populate the previously-empty "Case Sensitive" value (row 15, column B) with the
# text "true" (not the boolean) -- use a leading apostrophe so the engine stores it
# as a text shared string instead of auto-coercing to a boolean cell.
$refCell = $ws.Range("B14")
$cell = $ws.Range("B15")
$cell.Value = "'true"
$refCell.Copy()
$cell.PasteSpecial(-4122) # xlPasteFormats -- restore the original (unstyled) cell format
